$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1272.3077  # ALC H18
$ws.Cells.Item(18, 9).Value = 1140  # ALC I18
$ws.Cells.Item(18, 10).Value = 2000  # ALC J18
$ws.Cells.Item(18, 11).Value = 1140  # ALC K18
$ws.Cells.Item(18, 12).Value = 2000  # ALC L18
$ws.Cells.Item(18, 13).Value = -856  # ALC M18
$ws.Cells.Item(18, 14).Value = -2568  # ALC N18
$ws.Cells.Item(62, 8).Value = 37042036  # ALC H62
$ws.Cells.Item(62, 9).Value = 55560556  # ALC I62
$ws.Cells.Item(62, 10).Value = 5000  # ALC J62
$ws.Cells.Item(62, 11).Value = 55560556  # ALC K62
$ws.Cells.Item(62, 12).Value = 5000  # ALC L62
$ws.Cells.Item(62, 13).Value = -55559932  # ALC M62
$ws.Cells.Item(62, 14).Value = -6248  # ALC N62
$ws.Cells.Item(65, 8).Value = 37042036  # ALC H65
$ws.Cells.Item(65, 9).Value = 55560556  # ALC I65
$ws.Cells.Item(65, 10).Value = 5000  # ALC J65
$ws.Cells.Item(65, 11).Value = 277802780  # ALC K65
$ws.Cells.Item(65, 12).Value = 25000  # ALC L65
$ws.Cells.Item(65, 13).Value = -277799660  # ALC M65
$ws.Cells.Item(65, 14).Value = -31240  # ALC N65
$ws.Cells.Item(80, 8).Value = 1039.125  # ALC H80
$ws.Cells.Item(80, 9).Value = 1420.4  # ALC I80
$ws.Cells.Item(80, 10).Value = 865.8182  # ALC J80
$ws.Cells.Item(80, 11).Value = 4261.200000000001  # ALC K80
$ws.Cells.Item(80, 12).Value = 2597.4546  # ALC L80
$ws.Cells.Item(80, 13).Value = -3263.200000000001  # ALC M80
$ws.Cells.Item(80, 14).Value = -4593.4546  # ALC N80
$ws.Cells.Item(83, 8).Value = 1039.125  # ALC H83
$ws.Cells.Item(83, 9).Value = 1420.4  # ALC I83
$ws.Cells.Item(83, 10).Value = 865.8182  # ALC J83
$ws.Cells.Item(83, 11).Value = 12783.6  # ALC K83
$ws.Cells.Item(83, 12).Value = 7792.3638  # ALC L83
$ws.Cells.Item(83, 13).Value = -7791.6  # ALC M83
$ws.Cells.Item(83, 14).Value = -17776.3638  # ALC N83
$ws.Cells.Item(106, 8).Value = 7405.619  # ALC H106
$ws.Cells.Item(106, 9).Value = 7625.9  # ALC I106
$ws.Cells.Item(106, 11).Value = 7625.9  # ALC K106
$ws.Cells.Item(106, 13).Value = -6994.9  # ALC M106
$ws.Cells.Item(137, 8).Value = 1523.1957  # ALC H137
$ws.Cells.Item(137, 9).Value = 1432.8889  # ALC I137
$ws.Cells.Item(137, 10).Value = 1651.5264  # ALC J137
$ws.Cells.Item(137, 11).Value = 4298.6667  # ALC K137
$ws.Cells.Item(137, 12).Value = 4954.5792  # ALC L137
$ws.Cells.Item(137, 13).Value = -1748.6667  # ALC M137
$ws.Cells.Item(137, 14).Value = -10054.5792  # ALC N137
$ws.Cells.Item(138, 8).Value = 440815.44  # ALC H138
$ws.Cells.Item(138, 9).Value = 1194.12  # ALC I138
$ws.Cells.Item(138, 10).Value = 589336.1  # ALC J138
$ws.Cells.Item(138, 11).Value = 3582.36  # ALC K138
$ws.Cells.Item(138, 12).Value = 1768008.3  # ALC L138
$ws.Cells.Item(138, 13).Value = 1557.64  # ALC M138
$ws.Cells.Item(138, 14).Value = -1778288.3  # ALC N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 904  # ARM H10
$ws.Cells.Item(10, 9).Value = 904  # ARM I10
$ws.Cells.Item(10, 11).Value = 904  # ARM K10
$ws.Cells.Item(10, 13).Value = -734  # ARM M10
$ws.Cells.Item(32, 8).Value = 6251.4546  # ARM H32
$ws.Cells.Item(32, 9).Value = 5131.5713  # ARM I32
$ws.Cells.Item(32, 11).Value = 5131.5713  # ARM K32
$ws.Cells.Item(32, 13).Value = -4844.5713  # ARM M32

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 500  # BSM H12
$ws.Cells.Item(12, 9).Value = 500  # BSM I12
$ws.Cells.Item(12, 11).Value = 500  # BSM K12
$ws.Cells.Item(12, 13).Value = -332  # BSM M12
$ws.Cells.Item(20, 8).Value = 3455.7144  # BSM H20
$ws.Cells.Item(20, 9).Value = 3207.2727  # BSM I20
$ws.Cells.Item(20, 10).Value = 4366.6665  # BSM J20
$ws.Cells.Item(20, 11).Value = 3207.2727  # BSM K20
$ws.Cells.Item(20, 12).Value = 4366.6665  # BSM L20
$ws.Cells.Item(20, 13).Value = -2960.2727  # BSM M20
$ws.Cells.Item(20, 14).Value = -4860.6665  # BSM N20
$ws.Cells.Item(22, 8).Value = 398.91666  # BSM H22
$ws.Cells.Item(22, 10).Value = 438.7  # BSM J22
$ws.Cells.Item(22, 12).Value = 438.7  # BSM L22
$ws.Cells.Item(22, 14).Value = -784.7  # BSM N22
$ws.Cells.Item(88, 8).Value = 41999  # BSM H88
$ws.Cells.Item(88, 10).Value = 41999  # BSM J88
$ws.Cells.Item(88, 12).Value = 41999  # BSM L88
$ws.Cells.Item(88, 14).Value = -42811  # BSM N88
$ws.Cells.Item(91, 8).Value = 41999  # BSM H91
$ws.Cells.Item(91, 10).Value = 41999  # BSM J91
$ws.Cells.Item(91, 12).Value = 41999  # BSM L91
$ws.Cells.Item(91, 14).Value = -44807  # BSM N91
$ws.Cells.Item(103, 8).Value = 0  # BSM H103
$ws.Cells.Item(103, 10).Value = 0  # BSM J103
$ws.Cells.Item(103, 12).Value = 0  # BSM L103
$ws.Cells.Item(103, 14).ClearContents()  # BSM N103 (was -22172)
$ws.Cells.Item(105, 8).Value = 111112810  # BSM H105
$ws.Cells.Item(105, 9).Value = 125001650  # BSM I105
$ws.Cells.Item(105, 11).Value = 125001650  # BSM K105
$ws.Cells.Item(105, 13).Value = -124999903  # BSM M105
$ws.Cells.Item(134, 8).Value = 860.4483  # BSM H134
$ws.Cells.Item(134, 9).Value = 739.03705  # BSM I134
$ws.Cells.Item(134, 11).Value = 2217.11115  # BSM K134
$ws.Cells.Item(134, 13).Value = 317.8888499999998  # BSM M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1529.4419  # CRP H31
$ws.Cells.Item(31, 9).Value = 1446.8096  # CRP I31
$ws.Cells.Item(31, 10).Value = 5000  # CRP J31
$ws.Cells.Item(31, 11).Value = 1446.8096  # CRP K31
$ws.Cells.Item(31, 12).Value = 5000  # CRP L31
$ws.Cells.Item(31, 13).Value = -1151.8096  # CRP M31
$ws.Cells.Item(31, 14).Value = -5590  # CRP N31
$ws.Cells.Item(34, 8).Value = 1529.4419  # CRP H34
$ws.Cells.Item(34, 9).Value = 1446.8096  # CRP I34
$ws.Cells.Item(34, 10).Value = 5000  # CRP J34
$ws.Cells.Item(34, 11).Value = 1446.8096  # CRP K34
$ws.Cells.Item(34, 12).Value = 5000  # CRP L34
$ws.Cells.Item(34, 13).Value = -1244.8096  # CRP M34
$ws.Cells.Item(34, 14).Value = -5404  # CRP N34
$ws.Cells.Item(43, 8).Value = 14999.5  # CRP H43
$ws.Cells.Item(43, 10).Value = 14999.5  # CRP J43
$ws.Cells.Item(43, 12).Value = 14999.5  # CRP L43
$ws.Cells.Item(43, 14).Value = -15367.5  # CRP N43
$ws.Cells.Item(86, 8).Value = 4810205  # CRP H86
$ws.Cells.Item(86, 9).Value = 6698786.5  # CRP I86
$ws.Cells.Item(86, 10).Value = 88752  # CRP J86
$ws.Cells.Item(86, 11).Value = 6698786.5  # CRP K86
$ws.Cells.Item(86, 12).Value = 88752  # CRP L86
$ws.Cells.Item(86, 13).Value = -6697663.5  # CRP M86
$ws.Cells.Item(86, 14).Value = -90998  # CRP N86
$ws.Cells.Item(89, 8).Value = 4810205  # CRP H89
$ws.Cells.Item(89, 9).Value = 6698786.5  # CRP I89
$ws.Cells.Item(89, 10).Value = 88752  # CRP J89
$ws.Cells.Item(89, 11).Value = 33493932.5  # CRP K89
$ws.Cells.Item(89, 12).Value = 443760  # CRP L89
$ws.Cells.Item(89, 13).Value = -33488316.5  # CRP M89
$ws.Cells.Item(89, 14).Value = -454992  # CRP N89
$ws.Cells.Item(101, 8).Value = 14999.5  # CRP H101
$ws.Cells.Item(101, 10).Value = 14999.5  # CRP J101
$ws.Cells.Item(101, 12).Value = 14999.5  # CRP L101
$ws.Cells.Item(101, 14).Value = -21489.5  # CRP N101
$ws.Cells.Item(105, 8).Value = 1028.5  # CRP H105
$ws.Cells.Item(105, 9).Value = 1037.5  # CRP I105
$ws.Cells.Item(105, 10).Value = 1010.5  # CRP J105
$ws.Cells.Item(105, 11).Value = 1037.5  # CRP K105
$ws.Cells.Item(105, 12).Value = 1010.5  # CRP L105
$ws.Cells.Item(105, 13).Value = 709.5  # CRP M105
$ws.Cells.Item(105, 14).Value = -4504.5  # CRP N105
$ws.Cells.Item(122, 8).Value = 811.15  # CRP H122
$ws.Cells.Item(122, 9).Value = 756.8333  # CRP I122
$ws.Cells.Item(122, 10).Value = 1300  # CRP J122
$ws.Cells.Item(122, 11).Value = 2270.4999  # CRP K122
$ws.Cells.Item(122, 12).Value = 3900  # CRP L122
$ws.Cells.Item(122, 13).Value = 179.5001000000002  # CRP M122
$ws.Cells.Item(122, 14).Value = -8800  # CRP N122
$ws.Cells.Item(134, 8).Value = 11112162  # CRP H134
$ws.Cells.Item(134, 9).Value = 968.6486  # CRP I134
$ws.Cells.Item(134, 10).Value = 62501430  # CRP J134
$ws.Cells.Item(134, 11).Value = 2905.9458  # CRP K134
$ws.Cells.Item(134, 12).Value = 187504290  # CRP L134
$ws.Cells.Item(134, 13).Value = -370.9458  # CRP M134
$ws.Cells.Item(134, 14).Value = -187509360  # CRP N134
$ws.Cells.Item(141, 8).Value = 616463.9  # CRP H141
$ws.Cells.Item(141, 9).Value = 9999  # CRP I141
$ws.Cells.Item(141, 10).Value = 683848.9  # CRP J141
$ws.Cells.Item(141, 11).Value = 9999  # CRP K141
$ws.Cells.Item(141, 12).Value = 683848.9  # CRP L141
$ws.Cells.Item(141, 13).Value = -4819  # CRP M141
$ws.Cells.Item(141, 14).Value = -694208.9  # CRP N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1739.6842  # CUL H122
$ws.Cells.Item(122, 10).Value = 1887.75  # CUL J122
$ws.Cells.Item(122, 12).Value = 16989.75  # CUL L122
$ws.Cells.Item(122, 14).Value = -21889.75  # CUL N122
$ws.Cells.Item(131, 8).Value = 17858136  # CUL H131
$ws.Cells.Item(131, 9).Value = 142857740  # CUL I131
$ws.Cells.Item(131, 10).Value = 1049.9796  # CUL J131
$ws.Cells.Item(131, 11).Value = 428573220  # CUL K131
$ws.Cells.Item(131, 12).Value = 3149.9388  # CUL L131
$ws.Cells.Item(131, 13).Value = -428568180  # CUL M131
$ws.Cells.Item(131, 14).Value = -13229.9388  # CUL N131

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 9).Value = 845.5  # LTW I22
$ws.Cells.Item(22, 10).Value = 1500  # LTW J22
$ws.Cells.Item(22, 11).Value = 845.5  # LTW K22
$ws.Cells.Item(22, 12).Value = 1500  # LTW L22
$ws.Cells.Item(22, 13).Value = -550.5  # LTW M22
$ws.Cells.Item(22, 14).Value = -2090  # LTW N22
$ws.Cells.Item(27, 9).Value = 845.5  # LTW I27
$ws.Cells.Item(27, 10).Value = 1500  # LTW J27
$ws.Cells.Item(27, 11).Value = 845.5  # LTW K27
$ws.Cells.Item(27, 12).Value = 1500  # LTW L27
$ws.Cells.Item(27, 13).Value = -738.5  # LTW M27
$ws.Cells.Item(27, 14).Value = -1714  # LTW N27
$ws.Cells.Item(40, 8).Value = 3907.0476  # LTW H40
$ws.Cells.Item(40, 9).Value = 1979.9333  # LTW I40
$ws.Cells.Item(40, 11).Value = 1979.9333  # LTW K40
$ws.Cells.Item(40, 13).Value = -1843.9333  # LTW M40
$ws.Cells.Item(136, 8).Value = 1282.3448  # LTW H136
$ws.Cells.Item(136, 9).Value = 1153.4615  # LTW I136
$ws.Cells.Item(136, 10).Value = 2399.3333  # LTW J136
$ws.Cells.Item(136, 11).Value = 3460.3845  # LTW K136
$ws.Cells.Item(136, 12).Value = 7197.999899999999  # LTW L136
$ws.Cells.Item(136, 13).Value = -910.3844999999997  # LTW M136
$ws.Cells.Item(136, 14).Value = -12297.9999  # LTW N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 318.43478  # WVR H113
$ws.Cells.Item(113, 10).Value = 357.72726  # WVR J113
$ws.Cells.Item(113, 12).Value = 1073.18178  # WVR L113
$ws.Cells.Item(113, 14).Value = -5413.18178  # WVR N113
$ws.Cells.Item(135, 8).Value = 88343.2  # WVR H135
$ws.Cells.Item(135, 10).Value = 88343.2  # WVR J135
$ws.Cells.Item(135, 12).Value = 88343.2  # WVR L135
$ws.Cells.Item(135, 14).Value = -98483.2  # WVR N135

Write-Output "Applied all Kujata_Profits market data updates"